# Add a new reference entry (row 5 / sheet row 6) for a basic mathematical
# model example, together with a new "URL" column on the Reference Data
# table so the source link can be recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reference Data")
$lo = $ws.ListObjects.Item(1)

# --- fill in the new reference row ------------------------------------
$ws.Range("A6").Value = 5
$ws.Range("D6").Value = 2020
$ws.Range("F6").Value = "Modelling a Pandemic"
$ws.Range("G6").Value = "Article"

# --- add the new "URL" column to the table -----------------------------
$newCol = $lo.ListColumns.Add()
$ws.Range("J1").Value = "URL"

# Match the formatting used by the rest of the table's data columns
# (left/centre aligned, wrapped text) so the new column's cells line up
# with the others.
$ws.Range("J2:J49").HorizontalAlignment = -4131
$ws.Range("J2:J49").VerticalAlignment = -4108
$ws.Range("J2:J49").WrapText = $true

# Give the new column a sensible width, similar to the other wide columns.
$ws.Columns.Item(10).ColumnWidth = 32.53

$ws.Range("J6").Value = "https://towardsdatascience.com/modelling-a-pandemic-eb94025f248f"
$ws.Range("E6").Value = "Christian Graf"

# --- reflect the selection / active sheet at save time ------------------
$ws.Activate()
$ws.Range("D8").Select()
